$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 315 (shifts existing rows 315:387 down to 316:388)
$ws.Rows.Item(315).Insert()

# Populate the newly inserted row 315 with the weekly data.
# Categorical columns mirror what used to be in row 315 (now shifted to row 316).
$ws.Cells.Item(315, 1).Value = 3
$ws.Cells.Item(315, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(315, 3).Value = "Coquimbo"
$ws.Cells.Item(315, 4).Value = 44782
$ws.Cells.Item(315, 5).Value = 5
$ws.Cells.Item(315, 6).Value = 100112009
$ws.Cells.Item(315, 7).Value = "Acelga"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 230
$ws.Cells.Item(315, 11).Value = 3300
$ws.Cells.Item(315, 12).Value = 3500
$ws.Cells.Item(315, 13).Value = 3396
$ws.Cells.Item(315, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(315, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(315, 16).Value = 566
$ws.Cells.Item(315, 17).Value = 6
$ws.Cells.Item(315, 18).Value = "Hortaliza"
